$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 69 (shifts old rows 69-72 down to 70-73, and extends
# the SUM(F2:F69) / shared-formula references automatically)
$ws.Rows.Item(69).Insert()

# The existing last entry (row 68) had its end time corrected
$ws.Cells.Item(68, 5).Value = 0.80555555555555547

# Populate the newly inserted row 69 with the extra working-hours entry
$ws.Cells.Item(69, 1).Value = 2014
$ws.Cells.Item(69, 2).Value = 3
$ws.Cells.Item(69, 3).Value = 16
$ws.Cells.Item(69, 4).Value = 0.40625
$ws.Cells.Item(69, 5).Value = 0.54166666666666663
$ws.Cells.Item(69, 6).Formula = "=(E69-D69)*24*60"
$ws.Cells.Item(69, 7).Formula = "=F69/60"

# Restore the view/selection state recorded after the edit
$aw = $excel.ActiveWindow
$aw.ScrollRow = 43
$aw.ScrollColumn = 1
$ws.Range("F69").Select() | Out-Null
